$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AZ2").Value = 4440126.976
$ws.Range("BA2").Value = 4586783.744
$ws.Range("AZ3").Value = 1530573.056
$ws.Range("BA3").Value = 1978653.952
$ws.Range("AZ4").Value = 478404.992
$ws.Range("BA4").Value = 1209
$ws.Range("AZ5").Value = 110763
$ws.Range("BA5").Value = 854078.0159999999
$ws.Range("AZ6").Value = 833198.976
$ws.Range("BA6").Value = 1011971.008
$ws.Range("AZ7").Value = 18286
$ws.Range("BA7").Value = 32632
$ws.Range("AZ9").Value = 21085
$ws.Range("BA9").Value = 20576
$ws.Range("AZ10").Value = 4023
$ws.Range("BA10").Value = 1913
$ws.Range("AZ11").Value = 64812
$ws.Range("BA11").Value = 56275
$ws.Range("AZ12").Value = 2597945.088
$ws.Range("BA12").Value = 2481491.968
$ws.Range("AZ13").Value = 7
$ws.Range("BA13").Value = 7
$ws.Range("AZ16").Value = 2958
$ws.Range("BA16").Value = 1307
$ws.Range("AZ19").Value = 0
$ws.Range("BA19").Value = 0
$ws.Range("AZ21").Value = 205490
$ws.Range("BA21").Value = 0
$ws.Range("AZ22").Value = 195244.992
$ws.Range("BA22").Value = 16050
$ws.Range("AZ23").Value = 95427
$ws.Range("BA23").Value = 91584
$ws.Range("AZ24").Value = 20937
$ws.Range("BA24").Value = 19004
$ws.Range("AZ26").Value = 4440126.976
$ws.Range("BA26").Value = 4586783.744
$ws.Range("AZ27").Value = 361155.008
$ws.Range("BA27").Value = 1012030.016
$ws.Range("AZ28").Value = 57713
$ws.Range("BA28").Value = 50812
$ws.Range("AZ29").Value = 34244
$ws.Range("BA29").Value = 32904
$ws.Range("AZ30").Value = 36433
$ws.Range("BA30").Value = 55518
$ws.Range("AZ31").Value = 54557
$ws.Range("BA31").Value = 659761.024
$ws.Range("AZ33").Value = 0
$ws.Range("BA33").Value = 0
$ws.Range("AZ34").Value = 61666
$ws.Range("BA34").Value = 83060
$ws.Range("AZ35").Value = 116542
$ws.Range("BA35").Value = 129975
$ws.Range("AZ37").Value = 2157323.008
$ws.Range("BA37").Value = 1545632
$ws.Range("AZ38").Value = 627041.9840000001
$ws.Range("BA38").Value = 25008
$ws.Range("AZ39").Value = 0
$ws.Range("BA39").Value = 0
$ws.Range("AZ40").Value = 65558
$ws.Range("BA40").Value = 341155.008
$ws.Range("AZ41").Value = 148716
$ws.Range("BA41").Value = 284063.008
$ws.Range("AZ43").Value = 1316007.04
$ws.Range("BA43").Value = 895406.0159999999
$ws.Range("AZ46").Value = 15155
$ws.Range("BA46").Value = 6868
$ws.Range("AZ47").Value = 1906494.024
$ws.Range("BA47").Value = 2022254.048
$ws.Range("AZ48").Value = 981342.976
$ws.Range("BA48").Value = 981342.976
$ws.Range("AZ49").Value = 0
$ws.Range("BA49").Value = 0
$ws.Range("AZ51").Value = 1977380.992
$ws.Range("BA51").Value = 1642631.936
$ws.Range("AZ52").Value = 0
$ws.Range("BA52").Value = 0
$ws.Range("AZ53").Value = 0
$ws.Range("BA53").Value = 0
$ws.Range("AZ54").Value = 0
$ws.Range("BA54").Value = 0
$ws.Range("AZ55").Value = -1052230.016
$ws.Range("BA55").Value = -601721.024
$ws.Range("AZ59").Value = 182293.952
$ws.Range("BA59").Value = 283342.016
$ws.Range("AZ60").Value = -83529.992
$ws.Range("BA60").Value = -136631.008
$ws.Range("AZ61").Value = 98764.03200000001
$ws.Range("BA61").Value = 149710
$ws.Range("AZ62").Value = 521
$ws.Range("BA62").Value = -130
$ws.Range("AZ63").Value = 4652
$ws.Range("BA63").Value = -26445
$ws.Range("AZ65").Value = 33654.008
$ws.Range("BA65").Value = 2592
$ws.Range("AZ66").Value = -323606.016
$ws.Range("BA66").Value = -1742
$ws.Range("AZ67").Value = -15703
$ws.Range("BA67").Value = 590
$ws.Range("AZ68").Value = -5903
$ws.Range("BA68").Value = 110693
$ws.Range("AZ69").Value = 22969
$ws.Range("BA69").Value = 115204
$ws.Range("AZ70").Value = -28871.008
$ws.Range("BA70").Value = -4511
$ws.Range("AZ74").Value = -207620
$ws.Range("BA74").Value = 235268
$ws.Range("AZ75").Value = -33858
$ws.Range("BA75").Value = -30099
$ws.Range("AZ76").Value = -138568
$ws.Range("BA76").Value = -49479
$ws.Range("AZ79").Value = 9648
$ws.Range("BA79").Value = 760
$ws.Range("AZ80").Value = -333385.024
$ws.Range("BA80").Value = 156450
